$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = 7.306
$ws.Range("A9").Value = -21.74
$ws.Range("B9").Value = 5.88
$ws.Range("D9").Value = -7.852000000000001
$ws.Range("B11").Value = 5.816000000000001
$ws.Range("A13").Value = -22.08
$ws.Range("A16").Value = -21.253
$ws.Range("B16").Value = 5.417
$ws.Range("A18").Value = -21.767
$ws.Range("A20").Value = -20.308
$ws.Range("D22").Value = -7.943
$ws.Range("B23").Value = 7.650000000000001
$ws.Range("B24").Value = 5.238
$ws.Range("A26").Value = -21.321
$ws.Range("B26").Value = 5.973
$ws.Range("A27").Value = -21.667
$ws.Range("D27").Value = -7.924000000000001
$ws.Range("A29").Value = -21.346
$ws.Range("D29").Value = -7.385
$ws.Range("D32").Value = -7.322000000000001
$ws.Range("B34").Value = 7.226000000000001
$ws.Range("A35").Value = -21.618
$ws.Range("B35").Value = 6.092000000000001
$ws.Range("A36").Value = -20.927
$ws.Range("D37").Value = -7.598999999999999
$ws.Range("D38").Value = -8.092000000000002
$ws.Range("D39").Value = -7.406999999999999
$ws.Range("D41").Value = -8.191000000000001
$ws.Range("B44").Value = 5.641
$ws.Range("A45").Value = -21.238
$ws.Range("D45").Value = -8.004999999999999
$ws.Range("B48").Value = 5.63
$ws.Range("D48").Value = -7.736
$ws.Range("B49").Value = 5.986
$ws.Range("D51").Value = -8.273
$ws.Range("B52").Value = 4.975
$ws.Range("A55").Value = -22.207
$ws.Range("D56").Value = -8.090999999999999
$ws.Range("A57").Value = -21.91
$ws.Range("D57").Value = -7.927000000000001
$ws.Range("D61").Value = -7.910000000000001
$ws.Range("D64").Value = -7.75
$ws.Range("B66").Value = 4.878
$ws.Range("B67").Value = 5.139
$ws.Range("A69").Value = -21.291
$ws.Range("B73").Value = 6.299000000000001
$ws.Range("D75").Value = -8.012
$ws.Range("A76").Value = -20.392
$ws.Range("A78").Value = -21.015
$ws.Range("B78").Value = 7.157000000000001
$ws.Range("B80").Value = 8.300999999999998
$ws.Range("A82").Value = -21.718
$ws.Range("D82").Value = -8.187000000000001
$ws.Range("A83").Value = -21.509
$ws.Range("D90").Value = -7.081
$ws.Range("B91").Value = 5.509
$ws.Range("A93").Value = -21.453
$ws.Range("D93").Value = -7.037999999999999
$ws.Range("A97").Value = -21.271
$ws.Range("B97").Value = 5.396
$ws.Range("B99").Value = 4.925999999999999
$ws.Range("D102").Value = -7.957000000000001
$ws.Range("B104").Value = 7.273000000000001
$ws.Range("D105").Value = -7.934
